$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows 12-16, appended below the existing table (rows 1-11).
# Layout mirrors the existing rows: merged A/B cells for rows that share one
# "case study" (AirBnB -> rows 12-14), plain A/B cells for single-row cases
# (Google -> row 15, Grammarly -> row 16).
# ---------------------------------------------------------------------------

# --- Row 12 (AirBnB, part 1) -------------------------------------------------
$ws.Cells.Item(12,1).Value = "https://medium.com/airbnb-engineering/how-ai-text-generation-models-are-reshaping-customer-support-at-airbnb-a851db0b4fa3"

$ws.Cells.Item(12,2).Value = @"
AirBnB
проект автоматизации customer support
"@

$ws.Cells.Item(12,3).Value = @"
- подавали на вход запрос юзера и 1 статью из своей базы знаний
- промптом просили ответить, релевантна ли статья (т.е. только да\нет)
- в такой постановке зафайнтьюнили MT5
"@

$ws.Cells.Item(12,4).Value = @"
метрики для классификации, посчитатнные на тестовом датасете
+
АБ-тест работы на проде
"@

$ws.Cells.Item(12,5).Value = "взяли исторические данные о том как люди-саппортеры общались с клиентами"

$ws.Cells.Item(12,6).Value = "не было таких проблем, т.к. свели к классификации"

# --- Row 13 (AirBnB, part 2) -------------------------------------------------
$ws.Cells.Item(13,3).Value = @"
- отобрали типы вопросов, ответы на которые саппортеры выделяют в тексте обращения
- зафайнтьюнили модель для Question-Answer общения, но только для этих вопросов
- для каждого вопроса считали классификационные метрики
"@

$ws.Cells.Item(13,4).Value = "метрики для классификации, посчитатнные на тестовом датасете"

$ws.Cells.Item(13,5).Value = "взяли исторические данные о том как люди-саппортеры общались с клиентами"

$ws.Cells.Item(13,6).Value = "не было таких проблем, т.к. свели к классификации"

# --- Row 14 (AirBnB, part 3) -------------------------------------------------
$ws.Cells.Item(14,3).Value = @"
- из всей истории общения с клиентами на основе простой регулярки отобрали  те семплы, где саппортер пытыается перефразировать запрос клиента 
- кластризовали все парафразы, глазами просмотрели все кластеры и удалили те кластеры, которые содержали слишком общие и неполезные сообщения
- файнтьюнили T5 модель на парах (запрос клиента) - (парафраз от суппортера)
"@

$ws.Cells.Item(14,4).Value = @"
не говорят явно
скорее всего как-то субъективно
"@

$ws.Cells.Item(14,5).Value = "хитро фильтровали историческую выборку"

$ws.Cells.Item(14,6).Value = "заморочились с очисткой датасета"

# --- Row 15 (Google, grammar correction on mobile) --------------------------
$ws.Cells.Item(15,1).Value = "https://blog.research.google/2021/10/grammar-correction-as-you-type-on-pixel.html"

$ws.Cells.Item(15,2).Value = @"
Google
проект корректировки текста на мобилке
"@

$ws.Cells.Item(15,3).Value = @"
- наскраппили фраз из интернета
- прогнали их через большую модель для корректировки
- на получившемся датасете тренировали с 0 маленькую модель
"@

$ws.Cells.Item(15,4).Value = @"
не говорят явно
скорее всего сравнивали выходы маленькой и большой модели
"@

$ws.Cells.Item(15,5).Value = @"
- использовали публичные данные
- прогнали их через существующую модель
"@

$ws.Cells.Item(15,6).Value = "проблема с постепенным вводом текста: эвристически определяли, в какой момент разумно делать корректировку и показывать её пользователю"

# --- Row 16 (Grammarly, adversarial GEC) -------------------------------------
$ws.Cells.Item(16,1).Value = "https://www.grammarly.com/blog/engineering/adversarial-grammatical-error-correction/"

$ws.Cells.Item(16,2).Value = @"
Grammarly
генерация подсказок по редактированию текста
"@

$ws.Cells.Item(16,3).Value = @"
статья про то, как они свели поиск ошибок в тексте к работе GAN’а
нам не очень релевантно
"@

# ---------------------------------------------------------------------------
# Merge the "case" cells for the AirBnB rows (12-14), same pattern as the
# existing stitch-fix (2-3) / Microsoft (4-5) blocks.
# ---------------------------------------------------------------------------
$ws.Range("A12:A14").Merge()
$ws.Range("B12:B14").Merge()

# ---------------------------------------------------------------------------
# Formatting: the merged A/B cells for rows 12-14 are centered (horizontally
# and vertically); B additionally wraps text, matching the other merged
# blocks (A2:A3/B2:B3, A4:A5/B4:B5) already in the sheet.
# ---------------------------------------------------------------------------
$centerRange = $ws.Range("A12:B14")
$centerRange.HorizontalAlignment = -4108
$centerRange.VerticalAlignment = -4108

$ws.Range("A12:A14").WrapText = $false
$ws.Range("B12:B14").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights, matching the authored sizes for the wrapped text in each row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 57.45
$ws.Rows.Item(13).RowHeight = 57.45
$ws.Rows.Item(14).RowHeight = 91.5
$ws.Rows.Item(15).RowHeight = 68.65
$ws.Rows.Item(16).RowHeight = 58.2

# ---------------------------------------------------------------------------
# Misc sheet-level bookkeeping carried by the diff.
# ---------------------------------------------------------------------------
$ws.StandardWidth = 11.625
$ws.Range("C17").Select()
